$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.542.13"
$ws.Range("E2").Value = "  -1.90%  "
$ws.Range("D3").Value = "3.035.98"
$ws.Range("E3").Value = "  -1.85%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'554.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.82%  "
$ws.Range("D6").Value = "'140.61"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.30%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "3.035.22"
$ws.Range("E8").Value = "  -1.71%  "
$ws.Range("D9").Value = "'0.517"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.64%  "
$ws.Range("E10").Value = "  -0.27%  "
$ws.Range("D11").Value = "'6.15"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -13.61%  "
$ws.Range("E12").Value = "  +4.16%  "
$ws.Range("E13").Value = "  -0.65%  "
$ws.Range("D14").Value = "'35.29"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.06%  "
$ws.Range("D15").Value = "3.532.10"
$ws.Range("D16").Value = "63.557.72"
$ws.Range("E16").Value = "  -1.92%  "
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "'0.109"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.32%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.029.78"
$ws.Range("E18").Value = "  -2.24%  "
$ws.Range("E19").Value = "  -1.05%  "
$ws.Range("D20").Value = "'471.74"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.34%  "
$ws.Range("D21").Value = "'13.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.09%  "
$ws.Range("E22").Value = "  +0.43%  "
$ws.Range("D23").Value = "'14.46"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.91%  "
$ws.Range("E24").Value = "  -0.86%  "
$ws.Range("D25").Value = "'82.35"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.84%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").Value = "'2.77"
$ws.Range("D27").Style = "Normal"
$ws.Range("E28").Value = "  -1.89%  "
$ws.Range("E29").Value = "  -2.35%  "
$ws.Range("D30").Value = "'0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("E31").Value = "  -0.39%  "
$ws.Range("E32").Value = "  -1.34%  "
$ws.Range("E33").Value = "  -1.87%  "
$ws.Range("E34").Value = "  -1.09%  "
$ws.Range("D35").Value = "'6.15"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.11%  "
$ws.Range("D36").Value = "'54.66"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.22%  "
$ws.Range("E37").Value = "  -0.64%  "
$ws.Range("D38").Value = "'438.52"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.08%  "
$ws.Range("E39").Value = "  -2.05%  "
$ws.Range("D40").Value = "2.991.69"
$ws.Range("E40").Value = "  -0.47%  "
$ws.Range("E41").Value = "  -5.81%  "
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("D43").Value = "'8.23"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.35%  "
$ws.Range("D44").Value = "'0.267"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.44%  "
$ws.Range("D45").Value = "'27.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.20%  "
$ws.Range("E46").Value = "  +6.03%  "
$ws.Range("E48").Value = "  -0.18%  "
$ws.Range("D49").Value = "'118.08"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.43%  "
$ws.Range("D50").Value = "0.0₃0507"
$ws.Range("E50").Value = "  -1.77%  "
$ws.Range("E51").Value = "  -0.50%  "
